$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 622.3333
$ws.Range("I19").Value = 200
$ws.Range("J19").Value = 743
$ws.Range("K19").Value = 200
$ws.Range("L19").Value = 743
$ws.Range("M19").Value = -25
$ws.Range("N19").Value = -1093

$ws.Range("H28").Value = 221.72223
$ws.Range("I28").Value = 214.81482
$ws.Range("J28").Value = 242.44444
$ws.Range("K28").Value = 214.81482
$ws.Range("L28").Value = 242.44444
$ws.Range("M28").Value = 270.18518
$ws.Range("N28").Value = -1212.44444

$ws.Range("H33").Value = 152.48148
$ws.Range("I33").Value = 151.85
$ws.Range("J33").Value = 154.28572
$ws.Range("K33").Value = 151.85
$ws.Range("L33").Value = 154.28572
$ws.Range("M33").Value = 77.15000000000001
$ws.Range("N33").Value = -612.28572

$ws.Range("H88").Value = 3657.2666
$ws.Range("I88").Value = 3329.6667
$ws.Range("J88").Value = 3739.1667
$ws.Range("K88").Value = 3329.6667
$ws.Range("L88").Value = 3739.1667
$ws.Range("M88").Value = -2923.6667
$ws.Range("N88").Value = -4551.1667

$ws.Range("H91").Value = 3657.2666
$ws.Range("I91").Value = 3329.6667
$ws.Range("J91").Value = 3739.1667
$ws.Range("K91").Value = 3329.6667
$ws.Range("L91").Value = 3739.1667
$ws.Range("M91").Value = -1925.6667
$ws.Range("N91").Value = -6547.1667

$ws.Range("H99").Value = 861.7273
$ws.Range("I99").Value = 683.625
$ws.Range("J99").Value = 1336.6666
$ws.Range("K99").Value = 2050.875
$ws.Range("L99").Value = 4009.9998
$ws.Range("M99").Value = -552.875
$ws.Range("N99").Value = -7005.9998

$ws.Range("H100").Value = 3260.7896
$ws.Range("I100").Value = 2966.6667
$ws.Range("J100").Value = 3315.9375
$ws.Range("K100").Value = 2966.6667
$ws.Range("L100").Value = 3315.9375
$ws.Range("M100").Value = -2425.6667
$ws.Range("N100").Value = -4397.9375

$ws.Range("H101").Value = 901.4167
$ws.Range("I101").Value = 730.25
$ws.Range("J101").Value = 1243.75
$ws.Range("K101").Value = 2190.75
$ws.Range("L101").Value = 3731.25
$ws.Range("M101").Value = -568.75
$ws.Range("N101").Value = -6975.25

$ws.Range("H112").Value = 2286.25
$ws.Range("J112").Value = 2286.25
$ws.Range("L112").Value = 6858.75
$ws.Range("N112").Value = -9074.75

$ws.Range("H137").Value = 1416.8918
$ws.Range("I137").Value = 1183.2413
$ws.Range("J137").Value = 2263.875
$ws.Range("K137").Value = 3549.7239
$ws.Range("L137").Value = 6791.625
$ws.Range("M137").Value = -999.7239
$ws.Range("N137").Value = -11891.625


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 615
$ws.Range("I4").Value = 615
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 615
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -499
$ws.Range("N4").ClearContents()

$ws.Range("H32").Value = 12088.846
$ws.Range("I32").Value = 13328.848
$ws.Range("J32").Value = 2582.1667
$ws.Range("K32").Value = 13328.848
$ws.Range("L32").Value = 2582.1667
$ws.Range("M32").Value = -13041.848
$ws.Range("N32").Value = -3156.1667

$ws.Range("H74").Value = 1257.75
$ws.Range("I74").Value = 865.5
$ws.Range("K74").Value = 865.5
$ws.Range("M74").Value = 8.5

$ws.Range("H77").Value = 1257.75
$ws.Range("I77").Value = 865.5
$ws.Range("K77").Value = 4327.5
$ws.Range("M77").Value = 40.5

$ws.Range("H122").Value = 1941.8125
$ws.Range("I122").Value = 1790.6428
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 5371.928400000001
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -2921.928400000001
$ws.Range("N122").Value = -13900


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 39308
$ws.Range("J35").Value = 39308
$ws.Range("L35").Value = 39308
$ws.Range("N35").Value = -39928

$ws.Range("H64").Value = 181
$ws.Range("J64").Value = 180
$ws.Range("L64").Value = 180
$ws.Range("N64").Value = -630

$ws.Range("H67").Value = 181
$ws.Range("J67").Value = 180
$ws.Range("L67").Value = 180
$ws.Range("N67").Value = -1740

$ws.Range("H94").Value = 2064.818
$ws.Range("I94").Value = 619.75
$ws.Range("J94").Value = 3798.9
$ws.Range("K94").Value = 619.75
$ws.Range("L94").Value = 3798.9
$ws.Range("M94").Value = -168.75
$ws.Range("N94").Value = -4700.9

$ws.Range("H107").Value = 1651.75
$ws.Range("I107").Value = 1839.6875
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 1839.6875
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = 80.3125
$ws.Range("N107").Value = -4740

$ws.Range("H134").Value = 3321.3333
$ws.Range("I134").Value = 3385.6
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 10156.8
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -7621.799999999999
$ws.Range("N134").Value = -14070


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2508.5862
$ws.Range("I31").Value = 1662.9048
$ws.Range("J31").Value = 4728.5
$ws.Range("K31").Value = 1662.9048
$ws.Range("L31").Value = 4728.5
$ws.Range("M31").Value = -1367.9048
$ws.Range("N31").Value = -5318.5

$ws.Range("H34").Value = 2508.5862
$ws.Range("I34").Value = 1662.9048
$ws.Range("J34").Value = 4728.5
$ws.Range("K34").Value = 1662.9048
$ws.Range("L34").Value = 4728.5
$ws.Range("M34").Value = -1460.9048
$ws.Range("N34").Value = -5132.5

$ws.Range("H58").Value = 2056
$ws.Range("I58").Value = 2163.6
$ws.Range("J58").Value = 980
$ws.Range("K58").Value = 2163.6
$ws.Range("L58").Value = 980
$ws.Range("M58").Value = -1960.6
$ws.Range("N58").Value = -1386

$ws.Range("H69").Value = 23842
$ws.Range("I69").Value = 10763
$ws.Range("J69").Value = 50000
$ws.Range("K69").Value = 10763
$ws.Range("L69").Value = 50000
$ws.Range("M69").Value = -10014
$ws.Range("N69").Value = -51498

$ws.Range("H72").Value = 23842
$ws.Range("I72").Value = 10763
$ws.Range("J72").Value = 50000
$ws.Range("K72").Value = 32289
$ws.Range("L72").Value = 150000
$ws.Range("M72").Value = -28545
$ws.Range("N72").Value = -157488

$ws.Range("H122").Value = 3569
$ws.Range("I122").Value = 4087.7778
$ws.Range("J122").Value = 2902
$ws.Range("K122").Value = 12263.3334
$ws.Range("L122").Value = 8706
$ws.Range("M122").Value = -9813.3334
$ws.Range("N122").Value = -13606

$ws.Range("H134").Value = 2113.64
$ws.Range("I134").Value = 1964.75
$ws.Range("J134").Value = 2378.3333
$ws.Range("K134").Value = 5894.25
$ws.Range("L134").Value = 7134.999899999999
$ws.Range("M134").Value = -3359.25
$ws.Range("N134").Value = -12204.9999

$ws.Range("H136").Value = 2056
$ws.Range("I136").Value = 2163.6
$ws.Range("J136").Value = 980
$ws.Range("K136").Value = 6490.799999999999
$ws.Range("L136").Value = 2940
$ws.Range("M136").Value = -3940.799999999999
$ws.Range("N136").Value = -8040


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 14926667
$ws.Range("I131").Value = 228.33333
$ws.Range("J131").Value = 16394841
$ws.Range("K131").Value = 684.99999
$ws.Range("L131").Value = 49184523
$ws.Range("M131").Value = 4355.00001
$ws.Range("N131").Value = -49194603


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 487.625
$ws.Range("I107").Value = 230.1
$ws.Range("J107").Value = 671.5714
$ws.Range("K107").Value = 230.1
$ws.Range("L107").Value = 671.5714
$ws.Range("M107").Value = 1689.9
$ws.Range("N107").Value = -4511.5714


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1833.3334
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 1800
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 1800
$ws.Range("M46").Value = -1812
$ws.Range("N46").Value = -2176

$ws.Range("H122").Value = 9528838
$ws.Range("I122").Value = 4853
$ws.Range("J122").Value = 25005312
$ws.Range("K122").Value = 14559
$ws.Range("L122").Value = 75015936
$ws.Range("M122").Value = -12109
$ws.Range("N122").Value = -75020836


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 701.8889
$ws.Range("I107").Value = 670.7857
$ws.Range("J107").Value = 810.75
$ws.Range("K107").Value = 2012.3571
$ws.Range("L107").Value = 2432.25
$ws.Range("M107").Value = -92.35710000000017
$ws.Range("N107").Value = -6272.25
